$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty fields captured by the author's import:
#   H2 -> "Dirección de correspondencia"
#   J2 -> "Capítulo perteneciente"
#   L2 -> "Observaciones"
$ws.Range("H2").Value = "pruebad"
$ws.Range("J2").Value = "caracas"
$ws.Range("L2").Value = "ninguna"

# Keep the active cell/selection in sync with where the author ended up.
$ws.Range("L2").Select()
